$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (column F) for rows 3-5
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 172
$wsExhibit.Range("F4").Value = 759
$wsExhibit.Range("F5").Value = 65

# Sheet "全部类型" - update "想去人数" (column F) for rows 4-6
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 172
$wsAll.Range("F5").Value = 759
$wsAll.Range("F6").Value = 65
